$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 34) for 2025-09-18.
# Column A holds a date-like string that must stay plain text (matching the
# existing rows), so we build it via a text formula and then convert the
# formula to a static value to avoid Excel's automatic date recognition.
$ws.Range("A34").Formula = '="2025-09-18"'
$ws.Range("A34").Copy()
$ws.Range("A34").PasteSpecial(-4163)

$ws.Range("B34").Value = 59.08000183105469
$ws.Range("C34").Value = 711.2000122070312
$ws.Range("D34").Value = 337.8500061035156
